$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1759776536312849
$ws.Range("C2").Value = 0.5921787709497207
$ws.Range("J2").Value = 0.01675977653631285
$ws.Range("P2").Value = 0.1424581005586592
$ws.Range("S2").Value = 0.07262569832402235
$ws.Range("B3").Value = 0.004672897196261682
$ws.Range("C3").Value = 0.004672897196261682
$ws.Range("J3").Value = 0.03271028037383177
$ws.Range("P3").Value = 0.8130841121495327
$ws.Range("S3").Value = 0.1448598130841121
$ws.Range("J4").Value = 0.078125
$ws.Range("P4").Value = 0.78125
$ws.Range("S4").Value = 0.140625
$ws.Range("P5").Value = 0.8571428571428571
$ws.Range("S5").Value = 0.1428571428571428
$ws.Range("B6").Value = 0.05853658536585366
$ws.Range("D6").Value = 0.03902439024390244
$ws.Range("E6").Value = 0.004878048780487805
$ws.Range("F6").Value = 0.04878048780487805
$ws.Range("J6").Value = 0.2048780487804878
$ws.Range("O6").Value = 0.02926829268292683
$ws.Range("Q6").Value = 0.2439024390243902
$ws.Range("R6").Value = 0.06341463414634146
$ws.Range("S6").Value = 0.3073170731707317
$ws.Range("B7").Value = 0.1132075471698113
$ws.Range("D7").Value = 0.0330188679245283
$ws.Range("E7").Value = 0.004716981132075472
$ws.Range("F7").Value = 0.0660377358490566
$ws.Range("J7").Value = 0.1509433962264151
$ws.Range("O7").Value = 0.01886792452830189
$ws.Range("Q7").Value = 0.2028301886792453
$ws.Range("R7").Value = 0.08018867924528301
$ws.Range("S7").Value = 0.330188679245283
$ws.Range("B8").Value = 0.1075949367088608
$ws.Range("D8").Value = 0.02953586497890295
$ws.Range("E8").Value = 0.002109704641350211
$ws.Range("F8").Value = 0.0379746835443038
$ws.Range("J8").Value = 0.1118143459915612
$ws.Range("O8").Value = 0.02320675105485232
$ws.Range("Q8").Value = 0.2510548523206751
$ws.Range("R8").Value = 0.04852320675105485
$ws.Range("S8").Value = 0.3881856540084388
$ws.Range("B9").Value = 0.09642857142857143
$ws.Range("D9").Value = 0.025
$ws.Range("E9").Value = 0.003571428571428571
$ws.Range("F9").Value = 0.02857142857142857
$ws.Range("J9").Value = 0.1428571428571428
$ws.Range("O9").Value = 0.02142857142857143
$ws.Range("Q9").Value = 0.225
$ws.Range("R9").Value = 0.05
$ws.Range("S9").Value = 0.4071428571428571
$ws.Range("B10").Value = 0.1298421807747489
$ws.Range("D10").Value = 0.02223816355810617
$ws.Range("E10").Value = 0.002869440459110474
$ws.Range("F10").Value = 0.07317073170731707
$ws.Range("J10").Value = 0.1104734576757532
$ws.Range("O10").Value = 0.007173601147776184
$ws.Range("Q10").Value = 0.2209469153515065
$ws.Range("R10").Value = 0.07101865136298421
$ws.Range("S10").Value = 0.3622668579626973
$ws.Range("G11").Value = 0.1181556195965418
$ws.Range("J11").Value = 0.09798270893371758
$ws.Range("K11").Value = 0.1815561959654179
$ws.Range("L11").Value = 0.579250720461095
$ws.Range("S11").Value = 0.02305475504322766
$ws.Range("G12").Value = 0.6796116504854369
$ws.Range("J12").Value = 0.2572815533980582
$ws.Range("K12").Value = 0.009708737864077669
$ws.Range("L12").Value = 0.02427184466019417
$ws.Range("S12").Value = 0.02912621359223301
$ws.Range("G13").Value = 0.6346153846153846
$ws.Range("J13").Value = 0.3269230769230769
$ws.Range("S13").Value = 0.03846153846153846
$ws.Range("F15").Value = 0.004854368932038835
$ws.Range("H15").Value = 0.1116504854368932
$ws.Range("I15").Value = 0.0825242718446602
$ws.Range("J15").Value = 0.3689320388349515
$ws.Range("K15").Value = 0.1019417475728155
$ws.Range("M15").Value = 0.01456310679611651
$ws.Range("O15").Value = 0.03883495145631068
$ws.Range("S15").Value = 0.2766990291262136
$ws.Range("F16").Value = 0.007117437722419928
$ws.Range("H16").Value = 0.1921708185053381
$ws.Range("I16").Value = 0.09608540925266904
$ws.Range("J16").Value = 0.3309608540925267
$ws.Range("K16").Value = 0.1316725978647687
$ws.Range("M16").Value = 0.02135231316725979
$ws.Range("N16").Value = 0.007117437722419928
$ws.Range("O16").Value = 0.05338078291814947
$ws.Range("S16").Value = 0.1601423487544484
$ws.Range("F17").Value = 0.006932409012131715
$ws.Range("H17").Value = 0.145580589254766
$ws.Range("I17").Value = 0.1074523396880416
$ws.Range("J17").Value = 0.4246100519930676
$ws.Range("K17").Value = 0.1143847487001733
$ws.Range("M17").Value = 0.01906412478336222
$ws.Range("O17").Value = 0.03812824956672443
$ws.Range("S17").Value = 0.1438474870017331
$ws.Range("F18").Value = 0.01204819277108434
$ws.Range("H18").Value = 0.1867469879518072
$ws.Range("I18").Value = 0.1265060240963855
$ws.Range("J18").Value = 0.3915662650602409
$ws.Range("K18").Value = 0.0783132530120482
$ws.Range("M18").Value = 0.01807228915662651
$ws.Range("O18").Value = 0.06626506024096386
$ws.Range("S18").Value = 0.1204819277108434
$ws.Range("F19").Value = 0.01289398280802292
$ws.Range("H19").Value = 0.2041547277936963
$ws.Range("I19").Value = 0.1088825214899714
$ws.Range("J19").Value = 0.3524355300859599
$ws.Range("K19").Value = 0.1002865329512894
$ws.Range("M19").Value = 0.02363896848137536
$ws.Range("N19").Value = 0.001432664756446991
$ws.Range("O19").Value = 0.05802292263610315
$ws.Range("S19").Value = 0.1382521489971347
